$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.071.34"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.900.74"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5045"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3924"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09336"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "1.886.70"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.314"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001117"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.215"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").Value = "28.129.05"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.631"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("D27").Value = "2.105.17"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.086"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1066"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.643"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06641"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02422"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.272"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6387"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.995"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.273"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.027"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("E51").Value = "  -1.29%  "
